# Insert a new weekly record at row 59 ("Fruta / hortaliza, semanal").
# Excel's native row-insert shifts the existing rows 59..150 down to 60..151,
# preserving their values/styles untouched, then we populate the freshly
# inserted row 59 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(59).Insert()

$ws.Range("A59").Value2 = 11
$ws.Range("B59").Value2 = "Vega Monumental Concepción"
$ws.Range("C59").Value2 = "Bíobío"
$ws.Range("D59").Value2 = 44791
$ws.Range("E59").Value2 = 8
$ws.Range("F59").Value2 = 100112043
$ws.Range("G59").Value2 = "Pepino ensalada"
$ws.Range("H59").Value2 = "Sin especificar"
$ws.Range("I59").Value2 = "Primera"
$ws.Range("J59").Value2 = 180
$ws.Range("K59").Value2 = 22000
$ws.Range("L59").Value2 = 23000
$ws.Range("M59").Value2 = 22556
$ws.Range("N59").Value2 = "`$/caja 60 unidades"
$ws.Range("O59").Value2 = "Región de Arica y Parinacota"
$ws.Range("P59").Value2 = 376
$ws.Range("Q59").Value2 = 60
$ws.Range("R59").Value2 = "Hortaliza"
